$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = 7198
$ws.Range("C3").Value = 14400
$ws.Range("D3").Value = 21700
$ws.Range("E3").Value = 23700
$ws.Range("F3").Value = 26100
$ws.Range("G3").Value = 25900

$ws.Range("B8").Value = 15600
$ws.Range("C8").Value = 41000
$ws.Range("D8").Value = 81600
$ws.Range("E8").Value = 144000
$ws.Range("F8").Value = 180000
$ws.Range("G8").Value = 180000

$ws.Range("B13").Value = 4015
$ws.Range("C13").Value = 5128
$ws.Range("D13").Value = 5496
$ws.Range("E13").Value = 5955
$ws.Range("F13").Value = 6089
$ws.Range("G13").Value = 4773

$ws.Range("B18").Value = 152000
$ws.Range("C18").Value = 265000
$ws.Range("D18").Value = 368000
$ws.Range("E18").Value = 125000
$ws.Range("F18").Value = 139000
$ws.Range("G18").Value = 129000

$ws.Range("B23").Value = 3893
$ws.Range("C23").Value = 4641
$ws.Range("D23").Value = 12400
$ws.Range("E23").Value = 7434
$ws.Range("F23").Value = 11300
$ws.Range("G23").Value = 15200

$ws.Range("B28").Value = 144000
$ws.Range("C28").Value = 188000
$ws.Range("D28").Value = 394000
$ws.Range("E28").Value = 240000
$ws.Range("F28").Value = 689000
$ws.Range("G28").Value = 468000

$ws.Range("B33").Value = 2890
$ws.Range("C33").Value = 3130
$ws.Range("D33").Value = 3459
$ws.Range("E33").Value = 3750
$ws.Range("F33").Value = 3867
$ws.Range("G33").Value = 3880

$ws.Range("B38").Value = 58900
$ws.Range("C38").Value = 73500
$ws.Range("D38").Value = 78500
$ws.Range("E38").Value = 85800
$ws.Range("F38").Value = 99900
$ws.Range("G38").Value = 112000
